$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 101
$ws.Range("D2").Value = 85
$ws.Range("E2").Value = 0.8415841584158416
$ws.Range("F2").Value = 0.8333333333333334
$ws.Range("G2").Value = 0.09943492116856165
$ws.Range("H2").Value = 0.08286243430713471
$ws.Range("I2").Value = 453942.9050599152
$ws.Range("J2").Value = 164357.4527309576
$ws.Range("L2").Value = 164357.4527309576
$ws.Range("M2").Value = 618300.3577908728
$ws.Range("N2").Value = 10015935.0988
$ws.Range("O2").Value = 9608194.168699998
$ws.Range("P2").Value = 0.01640959641907515
$ws.Range("Q2").Value = 0.01710596703659199

# Row 3
$ws.Range("C3").Value = 103
$ws.Range("E3").Value = 0.8446601941747572
$ws.Range("G3").Value = 0.09776072146367178
$ws.Range("H3").Value = 0.08257458997416935
$ws.Range("I3").Value = 475461.5553898957
$ws.Range("J3").Value = 172518.3579019778
$ws.Range("L3").Value = 172518.3579019778
$ws.Range("M3").Value = 647979.9132918735
$ws.Range("N3").Value = 10570506.655464
$ws.Range("O3").Value = 10163133.497461
$ws.Range("P3").Value = 0.01632072742821664
$ws.Range("Q3").Value = 0.01697491801569636

# Row 4
$ws.Range("D4").Value = 88
$ws.Range("E4").Value = 0.8461538461538461
$ws.Range("F4").Value = 0.8461538461538461
$ws.Range("G4").Value = 0.09666630872252636
$ws.Range("H4").Value = 0.08179456891906077
$ws.Range("I4").Value = 499772.3434468232
$ws.Range("J4").Value = 177587.3793366524
$ws.Range("L4").Value = 177587.3793366524
$ws.Range("M4").Value = 677359.7227834756
$ws.Range("N4").Value = 10935308.90212792
$ws.Range("O4").Value = 10526964.54938483
$ws.Range("P4").Value = 0.01623981370129338
$ws.Range("Q4").Value = 0.01686976131662097

# Row 5
$ws.Range("D5").Value = 90
$ws.Range("E5").Value = 0.8571428571428571
$ws.Range("F5").Value = 0.8571428571428571
$ws.Range("G5").Value = 0.09525529745848933
$ws.Range("H5").Value = 0.0816473978215623
$ws.Range("I5").Value = 519331.4335515244
$ws.Range("J5").Value = 184792.9606174003
$ws.Range("L5").Value = 184792.9606174003
$ws.Range("M5").Value = 704124.3941689247
$ws.Range("N5").Value = 11301297.20929176
$ws.Range("O5").Value = 10890602.52596637
$ws.Range("P5").Value = 0.0163514822409472
$ws.Range("Q5").Value = 0.01696811174375338

# Row 6
$ws.Range("G6").Value = 0.0962254754163797
$ws.Range("H6").Value = 0.0798853003456737
$ws.Range("I6").Value = 530251.3219092456
$ws.Range("J6").Value = 187191.7221115101
$ws.Range("L6").Value = 187191.7221115101
$ws.Range("M6").Value = 717443.0440207556
$ws.Range("N6").Value = 11764042.88127051
$ws.Range("O6").Value = 11349577.35744537
$ws.Range("P6").Value = 0.01591219311258524
$ws.Range("Q6").Value = 0.01649327690503925
